$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest crypto
# quote data. Price values that look like plain numbers (e.g. "174.39") are
# written with a leading apostrophe so Excel keeps them as text, matching the
# source data (every Price/Volume cell in this sheet is stored as text, and
# several already contain thousand-separator dots like "62.969.69").
$ws.Range("D2").Value = "62.969.69"
$ws.Range("E2").Value = "  -4.36%  "
$ws.Range("D3").Value = "3.225.74"
$ws.Range("E3").Value = "  -5.29%  "
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "'174.39"
$ws.Range("E5").Value = "  -4.58%  "
$ws.Range("D6").Value = "'515.54"
$ws.Range("E6").Value = "  -3.66%  "
$ws.Range("D7").Value = "'0.591"
$ws.Range("E7").Value = "  -4.11%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "3.224.92"
$ws.Range("E9").Value = "  -5.18%  "
$ws.Range("D10").Value = "'0.599"
$ws.Range("E10").Value = "  -5.22%  "
$ws.Range("D11").Value = "'52.63"
$ws.Range("E11").Value = "  -9.22%  "
$ws.Range("D12").Value = "'0.130"
$ws.Range("E12").Value = "  -4.64%  "
$ws.Range("D13").Value = "'0.0000251"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").Value = "'8.87"
$ws.Range("E14").Value = "  -6.05%  "
$ws.Range("D15").Value = "3.734.84"
$ws.Range("E15").Value = "  -6.20%  "
$ws.Range("E16").Value = "  -6.31%  "
$ws.Range("D17").Value = "3.220.91"
$ws.Range("E17").Value = "  -6.25%  "
$ws.Range("D18").Value = "62.846.35"
$ws.Range("E18").Value = "  -4.60%  "
$ws.Range("D19").Value = "'17.14"
$ws.Range("E19").Value = "  -2.98%  "
$ws.Range("D20").Value = "'10.96"
$ws.Range("E20").Value = "  -3.90%  "
$ws.Range("D21").Value = "'0.956"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("D22").Value = "'366.02"
$ws.Range("E22").Value = "  -4.09%  "
$ws.Range("D23").Value = "'3.71"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").Value = "'79.97"
$ws.Range("E24").Value = "  -4.08%  "
$ws.Range("D25").Value = "'11.05"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "'3.89"
$ws.Range("E26").Value = "  +5.30%  "
$ws.Range("D27").Value = "'6.10"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").Value = "'2.61"
$ws.Range("E28").Value = "  -3.99%  "
$ws.Range("D29").Value = "'11.22"
$ws.Range("E29").Value = "  -5.00%  "
$ws.Range("D30").Value = "'8.15"
$ws.Range("E30").Value = "  -5.47%  "
$ws.Range("D31").Value = "'653.14"
$ws.Range("D32").Value = "'28.17"
$ws.Range("E32").Value = "  -5.91%  "
$ws.Range("D33").Value = "'6.29"
$ws.Range("E33").Value = "  -8.22%  "
$ws.Range("D34").Value = "'11.08"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").Value = "'0.104"
$ws.Range("E35").Value = "  -3.24%  "
$ws.Range("D36").Value = "'57.34"
$ws.Range("E36").Value = "  -7.21%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'36.50"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Value = "'0.373"
$ws.Range("E39").Value = "  -4.92%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").Value = "0.0₃0702"
$ws.Range("E41").Value = "  +11.60%  "
$ws.Range("D42").Value = "'0.122"
$ws.Range("E42").Value = "  -4.49%  "
$ws.Range("D43").Value = "2.860.21"
$ws.Range("E43").Value = "  -1.73%  "
$ws.Range("D44").Value = "'2.52"
$ws.Range("E44").Value = "  +4.93%  "
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").Value = "'2.86"
$ws.Range("E46").Value = "  +9.34%  "
$ws.Range("D47").Value = "'0.0390"
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").Value = "'2.59"
$ws.Range("E48").Value = "  -7.94%  "
$ws.Range("D49").Value = "'135.33"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("D50").Value = "'0.123"
$ws.Range("E50").Value = "  -2.83%  "
$ws.Range("D51").Value = "'2.94"
$ws.Range("E51").Value = "  +0.29%  "
